# Weekly fruit/vegetable price update:
# a new record (week of 2023-01-05) is inserted at row 175, pushing the
# existing rows 175-187 down to 176-188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 175, shifting rows 175:187 down to 176:188
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new weekly entry
$ws.Range("A175").Value = 2
$ws.Range("B175").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C175").Value = "Coquimbo"
$ws.Range("D175").Value = 44931
$ws.Range("E175").Value = 4
$ws.Range("F175").Value = 100112043
$ws.Range("G175").Value = "Pepino ensalada"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 700
$ws.Range("K175").Value = 11000
$ws.Range("L175").Value = 12000
$ws.Range("M175").Value = 11500
$ws.Range("N175").Value = "$/caja 70 unidades"
$ws.Range("O175").Value = "Provincia de Limarí"
$ws.Range("P175").Value = 164
$ws.Range("Q175").Value = 70
$ws.Range("R175").Value = "Hortaliza"
